$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds formatted price strings (e.g. "1.001", "30.700.48",
# "1.0000") that must stay text -- otherwise Excel coerces them to
# numbers and trailing/placeholder zeros or thousands-dot grouping is
# lost. Force Text format on each target cell individually before
# writing (a single multi-area "D2,D4,..." Range only honors the
# first area for NumberFormat, so each cell gets its own statement).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.700.48"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "241.88"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.4917"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("D8").Value = "0.2938"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").Value = "0.06750"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("D10").Value = "1.895.97"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").Value = "17.19"
$ws.Range("E11").Value = "  +5.22%  "
$ws.Range("D12").Value = "0.07261"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "90.96"
$ws.Range("E13").Value = "  +5.63%  "
$ws.Range("D14").Value = "0.6761"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("D15").Value = "5.036"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "30.679.62"
$ws.Range("E16").Value = "  +1.92%  "
$ws.Range("D17").Value = "0.000007997"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "13.15"
$ws.Range("E19").Value = "  +3.00%  "
$ws.Range("D20").Value = "2.140.43"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "4.805"
$ws.Range("E22").Value = "  +1.06%  "
$ws.Range("D23").Value = "191.20"
$ws.Range("E23").Value = "  +33.14%  "
$ws.Range("D24").Value = "6.098"
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").Value = "9.384"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "157.17"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").Value = "18.82"
$ws.Range("E27").Value = "  +11.05%  "
$ws.Range("E28").Value = "  +0.66%  "
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "4.298"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").Value = "0.09070"
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("D32").Value = "4.004"
$ws.Range("E32").Value = "  +0.75%  "
$ws.Range("D33").Value = "0.05243"
$ws.Range("E33").Value = "  +1.12%  "
$ws.Range("D34").Value = "0.7410"
$ws.Range("E34").Value = "  +3.01%  "
$ws.Range("D35").Value = "1.106"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").Value = "2.746"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("D37").Value = "0.01831"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "0.9339"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("D40").Value = "2.125"
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").Value = "0.4409"
$ws.Range("E41").Value = "  +3.73%  "
$ws.Range("D42").Value = "105.20"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").Value = "5.734"
$ws.Range("D45").Value = "0.1354"
$ws.Range("E45").Value = "  +5.47%  "
$ws.Range("D46").Value = "7.530"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").Value = "0.05873"
$ws.Range("E47").Value = "  +2.61%  "
$ws.Range("D48").Value = "8.775"
$ws.Range("E48").Value = "  +5.45%  "
$ws.Range("D49").Value = "1.431"
$ws.Range("E49").Value = "  +6.05%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").Value = "0.3958"
$ws.Range("E50").Value = "  +5.06%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "33.84"
$ws.Range("E51").Value = "  +2.97%  "
